$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = "-"
$ws.Range("C3").Value  = "[-, -, -, 'MEC-3B-M.A.Comp.CAD / CAM']"

$ws.Range("B4").Value  = "-"
$ws.Range("C4").Value  = "[-, -, -, 'MEC-3B-M.A.Comp.CAD / CAM']"

$ws.Range("B6").Value  = "-"
$ws.Range("C6").Value  = "[-, -, -, 'MEC-3B-M.A.Comp.CAD / CAM']"

$ws.Range("B7").Value  = "-"
$ws.Range("C7").Value  = "[-, -, -, 'MEC-3B-M.A.Comp.CAD / CAM']"

$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "[-, 'MEC-3A-Usin. CNC', 'MEC-3A-M.A.Comp.CAD / CAM', -]"

$ws.Range("D11").Value = "['MEC-3A-M.A.Comp.CAD / CAM', -, 'MEC-3A-Usin. CNC', -]"
$ws.Range("E11").Value = "-"

$ws.Range("E12").Value = "-"

$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "-"

$ws.Range("E15").Value = "-"

$ws.Range("E16").Value = "['MEC-3A-M.A.Comp.CAD / CAM', -, -, -]"
$ws.Range("F16").Value = "['MEC-3A-M.A.Comp.CAD / CAM', -, -, 'MEC-3A-Usin. CNC']"

$ws.Range("B18").Value = "[-, -, -, 'ELM-2NA-CAM']"
$ws.Range("C18").Value = "['ELM-2NA-CAM', -, -, -]"
$ws.Range("E18").Value = "['MEC-2NB-CAD/CAM', -, -, -]"
$ws.Range("F18").Value = "['MEC-2NA-Usin. CNC', -, -, -]"

$ws.Range("B19").Value = "[-, -, -, 'ELM-2NA-CAM']"
$ws.Range("E19").Value = "['MEC-2NB-CAD/CAM', -, -, -]"
$ws.Range("F19").Value = "['MEC-2NA-Usin. CNC', -, -, -]"

$ws.Range("B20").Value = "[-, -, -, 'ELM-2NA-CAM']"
$ws.Range("E20").Value = "['MEC-2NB-CAD/CAM', -, -, -]"
$ws.Range("F20").Value = "['MEC-2NA-Usin. CNC', -, -, -]"

$ws.Range("B21").Value = "-"
$ws.Range("E21").Value = "[-, 'MEC-2NB-CAD/CAM', -, -]"
$ws.Range("F21").Value = "['MEC-2NA-Usin. CNC', -, -, -]"
